$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("圣邦股份", "圣邦股份", "圣邦股份"),
    @("卧龙电驱", "上海贝岭", "卧龙电驱"),
    @("上海贝岭", "上海建工", "万通发展"),
    @("首开股份", "卧龙电驱", "上海建工"),
    @("上海建工", "山子高科", "吉视传媒"),
    @("方正科技", "北方铜业", "山子高科"),
    @("青山纸业", "方正科技", "首开股份"),
    @("吉视传媒", "思瑞浦", "方正科技"),
    @("山子高科", "北方稀土", "青山纸业"),
    @("思瑞浦", "北京君正", "金发科技"),
    @("北方稀土", "吉视传媒", "北方稀土"),
    @("北京君正", "大洋电机", "利欧股份"),
    @("北方铜业", "兆易创新", "华胜天成"),
    @("兆易创新", "首开股份", "兆驰股份"),
    @("万通发展", "工业富联", "工业富联"),
    @("工业富联", "东方财富", "上海贝岭"),
    @("香农芯创", "士兰微", "立讯精密"),
    @("大洋电机", "金发科技", "先导智能"),
    @("金发科技", "青山纸业", "胜宏科技"),
    @("先导智能", "兆驰股份", "国轩高科")
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    $ws.Cells.Item($row, 1).Value = $rowVals[0]
    $ws.Cells.Item($row, 2).Value = $rowVals[1]
    $ws.Cells.Item($row, 3).Value = $rowVals[2]
}
